$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the "总计" sheet,
#    mirroring the layout used by the "2021-Q1" sheet.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")   # "总计" - new sheet goes right before it

$new = $wb.Worksheets.Add($wsTotal)
$new.Name = "2022-Q1"

# NOTE: worksheet object references resolve by their original position,
# so after inserting/renaming sheets we must look sheets back up by name
# (rather than reuse older handles) before touching them again.
$new = $wb.Worksheets.Item("2022-Q1")

$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

$new.Range("A2").Value = 0

# These look numeric but must be stored as text, matching the source data.
$new.Range("B2:G2").NumberFormat = "@"
$new.Range("B2").Value = "202801"
$new.Range("C2").Value = "南方全球精选配置(QDII-FOF)"
$new.Range("D2").Value = "18.00"
$new.Range("E2").Value = "28.82"
$new.Range("F2").Value = "0.93"
$new.Range("G2").Value = "0.1674"
$new.Range("H2").Value = 10
# Drop the temporary "@" text format so these cells end up unstyled, like the source.
$new.Range("B2:G2").ClearFormats()

# Copy the header / index-column formatting (style used on "2021-Q1") onto the new sheet.
$wsQ1 = $wb.Worksheets.Item("2021-Q1")
$wsQ1.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)
$wsQ1.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing rows down and renumbering the leading index column.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.17

# Restore the index-column style on the newly inserted row and renumber
# the existing rows that shifted down.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

Write-Host "done"
